$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New explanation text for column B, rows 10-66 (per target codebook revision)
$explanations = @{
    10 = "[excluded vs. included] If the test subject did not have 120 in `"total_trials`", then it was excluded."
    11 = "This variable notes the average of how many of the trials are under 300 ms."
    12 = "Notes the average number of times a subject has given the correct answer in their trial."
    13 = "[excluded vs. included] This variable excludes the subjects who had more than 10% of the trials < 300ms."
    14 = "[excluded vs. included]This variable excludes subjects who had less than 75% accurate in their trials."
    15 = "[1-6] Item from the Agreeableness Scale"
    16 = "[1-6] Item from the Agreeableness Scale"
    17 = "[1-6] Item from the Agreeableness Scale"
    18 = "[1-6] Item from the Agreeableness Scale"
    19 = "[1-6] Item from the Agreeableness Scale"
    20 = "[1-6] Item from the Agreeableness Scale"
    21 = "[1-6] Item from the Agreeableness Scale"
    22 = "[1-6] Item from the Agreeableness Scale"
    23 = "[1-6] Item from the Agreeableness Scale"
    24 = "[1-6] Item from the Openness Scale"
    25 = "[1-6] Item from the Openness Scale"
    26 = "[1-6] Item from the Openness Scale"
    27 = "[1-6] Item from the Openness Scale"
    28 = "[1-6] Item from the Openness Scale"
    29 = "[1-6] Item from the Openness Scale"
    30 = "[1-6] Item from the Openness Scale"
    31 = "[1-6] Item from the Openness Scale"
    32 = "[1-6] Item from the Openness Scale"
    33 = "[1-6] Item from the Openness Scale"
    34 = "[1-6] Item from the Conscientiousness Scale"
    35 = "[1-6] Item aus der Conscientiousness Skala"
    36 = "[1-6] Item aus der Conscientiousness Skala"
    37 = "[1-6] Item aus der Conscientiousness Skala"
    38 = "[1-6] Item aus der Conscientiousness Skala"
    39 = "[1-6] Item aus der Conscientiousness Skala"
    40 = "[1-6] Item aus der Conscientiousness Skala"
    41 = "[1-6] Item aus der Conscientiousness Skala"
    42 = "[1-6] Item aus der Conscientiousness Skala"
    43 = "[1-6] Item from the Extraversion Scale"
    44 = "[1-6] Item from the Extraversion Scale"
    45 = "[1-6] Item from the Extraversion Scale"
    46 = "[1-6] Item from the Extraversion Scale"
    47 = "[1-6] Item from the Extraversion Scale"
    48 = "[1-6] Item from the Extraversion Scale"
    49 = "[1-6] Item from the Extraversion Scale"
    50 = "[1-6] Item from the Extraversion Scale"
    51 = "[1-6] Item from the Neuroticism Scale"
    52 = "[1-6] Item from the Neuroticism Scale"
    53 = "[1-6] Item from the Neuroticism Scale"
    54 = "[1-6] Item from the Neuroticism Scale"
    55 = "[1-6] Item from the Neuroticism Scale"
    56 = "[1-6] Item from the Neuroticism Scale"
    57 = "[1-6] Item from the Neuroticism Scale"
    58 = "[1-6] Item from the Neuroticism Scale"
    59 = "[excluded vs. included] This variable excludes subjects who have given an impossible answer (e.g. 1-6 are possible, but 7 or 0 are impossible)."
    60 = "[excluded vs. included] This variable excludes subjects who have incompletely filled out a scale (note: it does not exclude them if they have not filled out an entire scale, but only if a scale is partially but not completely filled out!)."
    61 = "The average score of a subject who completed the Agreeableness Scale. (all bfi_a items/ by the number of items)."
    62 = "The average score of a subject who completed the Openness Scale.  (all bfi_o items/ divided by the number of items)."
    63 = "The average score of a subject who completed the Extraversion Scale.  (all bfi_e items/ by the number of items)."
    64 = "The average score of a subject who completed the Neuroticism Scale.  (all bfi_n items/ by the number of items)"
    65 = "The average score of a subject who completed the Conscientiousness Scale.  (all bfi_c items/ divided by the number of items)."
    66 = "The average score of a subject who completed the Agreeableness Scale. (all bfi_a items/ by the number of items)."
}

foreach ($row in $explanations.Keys) {
    $ws.Cells.Item($row, 2).Value = $explanations[$row]
}

# B35 previously carried a redundant non-bold style (cellXfs index 2); normalize it
# back to the default/general style so it matches its sibling rows (e.g. B36).
$ws.Range("B35").ClearFormats()
$ws.Range("B35").Value = $explanations["35"]

# Restore the text font/size etc. is inherited from defaults already; re-apply value

# Update the active selection to match the saved view (B10).
$ws.Range("B10").Select()
